$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: "100nF" debounce capacitors - add C52 C53 C54 C55 to the reference
# designator list and bump the quantity from 33 to 37 (33 + 4 new caps).
$ws.Range("A4").Value = "C3 C4 C5 C6 C8 C9 C14 C15 C16 C17 C18 C19 C20 C21 C22 C23 C24 C25 C27 C28 C29 C30 C32 C33 C34 C35 C36 C37 C38 C39 C40 C41 C42 C52 C53 C54 C55 "
$ws.Range("B4").Value = 37

# Row 23: trigger-input resistor group value/part changed from 1K (C17513)
# to 2.2K (C17520) to pair with the new debounce caps.
$ws.Range("C23").Value = "2.2K"
$ws.Range("D23").Value = "C17520"
